$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-05 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-06 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("20+11=", $true, $false, $false, $false, $false, $true, 1, $false, "94-44=", 2) | Out-Null
$d.Content.Find.Execute("90-61=", $true, $false, $false, $false, $false, $true, 1, $false, "63-13=", 2) | Out-Null
$d.Content.Find.Execute("78-71=", $true, $false, $false, $false, $false, $true, 1, $false, "73+25=", 2) | Out-Null
$d.Content.Find.Execute("92-11=", $true, $false, $false, $false, $false, $true, 1, $false, "47+15=", 2) | Out-Null
$d.Content.Find.Execute("62-7=", $true, $false, $false, $false, $false, $true, 1, $false, "66-8=", 2) | Out-Null
$d.Content.Find.Execute("36-9=", $true, $false, $false, $false, $false, $true, 1, $false, "32+56=", 2) | Out-Null
$d.Content.Find.Execute("99-64=", $true, $false, $false, $false, $false, $true, 1, $false, "91-51=", 2) | Out-Null
$d.Content.Find.Execute("88-10=", $true, $false, $false, $false, $false, $true, 1, $false, "65+23=", 2) | Out-Null
$d.Content.Find.Execute("90+4=", $true, $false, $false, $false, $false, $true, 1, $false, "95-43=", 2) | Out-Null
$d.Content.Find.Execute("23+11=", $true, $false, $false, $false, $false, $true, 1, $false, "93-75=", 2) | Out-Null
$d.Content.Find.Execute("54+34=", $true, $false, $false, $false, $false, $true, 1, $false, "66-18=", 2) | Out-Null
$d.Content.Find.Execute("3+91=", $true, $false, $false, $false, $false, $true, 1, $false, "79-78=", 2) | Out-Null
$d.Content.Find.Execute("99-86=", $true, $false, $false, $false, $false, $true, 1, $false, "94-15=", 2) | Out-Null
$d.Content.Find.Execute("31+7=", $true, $false, $false, $false, $false, $true, 1, $false, "9+35=", 2) | Out-Null
$d.Content.Find.Execute("48-29=", $true, $false, $false, $false, $false, $true, 1, $false, "15+30=", 2) | Out-Null
$d.Content.Find.Execute("68-13=", $true, $false, $false, $false, $false, $true, 1, $false, "82+4=", 2) | Out-Null
$d.Content.Find.Execute("3+73=", $true, $false, $false, $false, $false, $true, 1, $false, "54+19=", 2) | Out-Null
$d.Content.Find.Execute("82+7=", $true, $false, $false, $false, $false, $true, 1, $false, "43-11=", 2) | Out-Null
$d.Content.Find.Execute("78+12=", $true, $false, $false, $false, $false, $true, 1, $false, "37-35=", 2) | Out-Null
$d.Content.Find.Execute("89+10=", $true, $false, $false, $false, $false, $true, 1, $false, "17+52=", 2) | Out-Null
$d.Content.Find.Execute("62-28=", $true, $false, $false, $false, $false, $true, 1, $false, "5+77=", 2) | Out-Null
$d.Content.Find.Execute("78-10=", $true, $false, $false, $false, $false, $true, 1, $false, "37+27=", 2) | Out-Null
$d.Content.Find.Execute("31-30=", $true, $false, $false, $false, $false, $true, 1, $false, "27+36=", 2) | Out-Null
$d.Content.Find.Execute("67-49=", $true, $false, $false, $false, $false, $true, 1, $false, "43+44=", 2) | Out-Null
$d.Content.Find.Execute("54-30=", $true, $false, $false, $false, $false, $true, 1, $false, "49-2=", 2) | Out-Null
$d.Content.Find.Execute("73-27=", $true, $false, $false, $false, $false, $true, 1, $false, "18-16=", 2) | Out-Null
$d.Content.Find.Execute("67-59=", $true, $false, $false, $false, $false, $true, 1, $false, "70-44=", 2) | Out-Null
$d.Content.Find.Execute("58+1=", $true, $false, $false, $false, $false, $true, 1, $false, "98-65=", 2) | Out-Null
$d.Content.Find.Execute("19+78=", $true, $false, $false, $false, $false, $true, 1, $false, "94-87=", 2) | Out-Null
$d.Content.Find.Execute("91-47=", $true, $false, $false, $false, $false, $true, 1, $false, "89-51=", 2) | Out-Null
$d.Content.Find.Execute("1+33=", $true, $false, $false, $false, $false, $true, 1, $false, "77-6=", 2) | Out-Null
$d.Content.Find.Execute("82-50=", $true, $false, $false, $false, $false, $true, 1, $false, "57-36=", 2) | Out-Null
$d.Content.Find.Execute("69-59=", $true, $false, $false, $false, $false, $true, 1, $false, "31+38=", 2) | Out-Null
$d.Content.Find.Execute("60+9=", $true, $false, $false, $false, $false, $true, 1, $false, "47+0=", 2) | Out-Null
$d.Content.Find.Execute("28+23=", $true, $false, $false, $false, $false, $true, 1, $false, "18-6=", 2) | Out-Null
$d.Content.Find.Execute("82-6=", $true, $false, $false, $false, $false, $true, 1, $false, "74-64=", 2) | Out-Null
$d.Content.Find.Execute("20+27=", $true, $false, $false, $false, $false, $true, 1, $false, "96-31=", 2) | Out-Null
$d.Content.Find.Execute("55-19=", $true, $false, $false, $false, $false, $true, 1, $false, "24+56=", 2) | Out-Null
$d.Content.Find.Execute("87-40=", $true, $false, $false, $false, $false, $true, 1, $false, "23+72=", 2) | Out-Null
$d.Content.Find.Execute("45+54=", $true, $false, $false, $false, $false, $true, 1, $false, "68-27=", 2) | Out-Null
$d.Content.Find.Execute("25+30=", $true, $false, $false, $false, $false, $true, 1, $false, "9+44=", 2) | Out-Null
$d.Content.Find.Execute("0+82=", $true, $false, $false, $false, $false, $true, 1, $false, "22-21=", 2) | Out-Null
$d.Content.Find.Execute("86-0=", $true, $false, $false, $false, $false, $true, 1, $false, "14+0=", 2) | Out-Null
$d.Content.Find.Execute("8+18=", $true, $false, $false, $false, $false, $true, 1, $false, "0+32=", 2) | Out-Null
$d.Content.Find.Execute("51-38=", $true, $false, $false, $false, $false, $true, 1, $false, "60+12=", 2) | Out-Null
$d.Content.Find.Execute("72-21=", $true, $false, $false, $false, $false, $true, 1, $false, "25-19=", 2) | Out-Null
$d.Content.Find.Execute("84-77=", $true, $false, $false, $false, $false, $true, 1, $false, "49+4=", 2) | Out-Null
$d.Content.Find.Execute("93-37=", $true, $false, $false, $false, $false, $true, 1, $false, "70-7=", 2) | Out-Null
$d.Content.Find.Execute("36-27=", $true, $false, $false, $false, $false, $true, 1, $false, "17+60=", 2) | Out-Null
$d.Content.Find.Execute("16-14=", $true, $false, $false, $false, $false, $true, 1, $false, "47-12=", 2) | Out-Null
$d.Content.Find.Execute("30+20=", $true, $false, $false, $false, $false, $true, 1, $false, "79-6=", 2) | Out-Null
$d.Content.Find.Execute("69-10=", $true, $false, $false, $false, $false, $true, 1, $false, "58+16=", 2) | Out-Null
$d.Content.Find.Execute("22-18=", $true, $false, $false, $false, $false, $true, 1, $false, "24+69=", 2) | Out-Null
$d.Content.Find.Execute("8+35=", $true, $false, $false, $false, $false, $true, 1, $false, "76+11=", 2) | Out-Null
$d.Content.Find.Execute("41+14=", $true, $false, $false, $false, $false, $true, 1, $false, "62-59=", 2) | Out-Null
$d.Content.Find.Execute("11+45=", $true, $false, $false, $false, $false, $true, 1, $false, "85+7=", 2) | Out-Null
$d.Content.Find.Execute("63-22=", $true, $false, $false, $false, $false, $true, 1, $false, "48-46=", 2) | Out-Null
$d.Content.Find.Execute("5+92=", $true, $false, $false, $false, $false, $true, 1, $false, "0+81=", 2) | Out-Null
$d.Content.Find.Execute("77-25=", $true, $false, $false, $false, $false, $true, 1, $false, "46-5=", 2) | Out-Null
$d.Content.Find.Execute("93-47=", $true, $false, $false, $false, $false, $true, 1, $false, "18+18=", 2) | Out-Null
$d.Content.Find.Execute("44+18=", $true, $false, $false, $false, $false, $true, 1, $false, "0+52=", 2) | Out-Null
$d.Content.Find.Execute("45+52=", $true, $false, $false, $false, $false, $true, 1, $false, "41-20=", 2) | Out-Null
$d.Content.Find.Execute("19-15=", $true, $false, $false, $false, $false, $true, 1, $false, "79-49=", 2) | Out-Null
$d.Content.Find.Execute("93-4=", $true, $false, $false, $false, $false, $true, 1, $false, "11+58=", 2) | Out-Null
$d.Content.Find.Execute("57-19=", $true, $false, $false, $false, $false, $true, 1, $false, "17+21=", 2) | Out-Null
$d.Content.Find.Execute("62-47=", $true, $false, $false, $false, $false, $true, 1, $false, "39+35=", 2) | Out-Null
$d.Content.Find.Execute("57+11=", $true, $false, $false, $false, $false, $true, 1, $false, "17+72=", 2) | Out-Null
$d.Content.Find.Execute("72-20=", $true, $false, $false, $false, $false, $true, 1, $false, "3+49=", 2) | Out-Null
$d.Content.Find.Execute("81-62=", $true, $false, $false, $false, $false, $true, 1, $false, "74+11=", 2) | Out-Null
$d.Content.Find.Execute("25+34=", $true, $false, $false, $false, $false, $true, 1, $false, "9+39=", 2) | Out-Null
$d.Content.Find.Execute("10+71=", $true, $false, $false, $false, $false, $true, 1, $false, "71-0=", 2) | Out-Null
$d.Content.Find.Execute("57-39=", $true, $false, $false, $false, $false, $true, 1, $false, "56+0=", 2) | Out-Null
$d.Content.Find.Execute("37+36=", $true, $false, $false, $false, $false, $true, 1, $false, "43+18=", 2) | Out-Null
$d.Content.Find.Execute("64-57=", $true, $false, $false, $false, $false, $true, 1, $false, "82+16=", 2) | Out-Null
$d.Content.Find.Execute("17+2=", $true, $false, $false, $false, $false, $true, 1, $false, "3+26=", 2) | Out-Null
$d.Content.Find.Execute("34-1=", $true, $false, $false, $false, $false, $true, 1, $false, "79-73=", 2) | Out-Null
$d.Content.Find.Execute("25+62=", $true, $false, $false, $false, $false, $true, 1, $false, "25+22=", 2) | Out-Null
$d.Content.Find.Execute("19-12=", $true, $false, $false, $false, $false, $true, 1, $false, "37+25=", 2) | Out-Null
$d.Content.Find.Execute("38+9=", $true, $false, $false, $false, $false, $true, 1, $false, "87-6=", 2) | Out-Null
$d.Content.Find.Execute("91-75=", $true, $false, $false, $false, $false, $true, 1, $false, "25-0=", 2) | Out-Null
$d.Content.Find.Execute("56-1=", $true, $false, $false, $false, $false, $true, 1, $false, "14+22=", 2) | Out-Null
$d.Content.Find.Execute("91+2=", $true, $false, $false, $false, $false, $true, 1, $false, "37-15=", 2) | Out-Null
$d.Content.Find.Execute("14+71=", $true, $false, $false, $false, $false, $true, 1, $false, "10+8=", 2) | Out-Null
$d.Content.Find.Execute("70-69=", $true, $false, $false, $false, $false, $true, 1, $false, "63+13=", 2) | Out-Null
$d.Content.Find.Execute("33+49=", $true, $false, $false, $false, $false, $true, 1, $false, "98-26=", 2) | Out-Null
$d.Content.Find.Execute("11+41=", $true, $false, $false, $false, $false, $true, 1, $false, "38+37=", 2) | Out-Null
$d.Content.Find.Execute("55+12=", $true, $false, $false, $false, $false, $true, 1, $false, "43+49=", 2) | Out-Null
$d.Content.Find.Execute("7+38=", $true, $false, $false, $false, $false, $true, 1, $false, "37+59=", 2) | Out-Null
$d.Content.Find.Execute("19+12=", $true, $false, $false, $false, $false, $true, 1, $false, "1+24=", 2) | Out-Null
$d.Content.Find.Execute("94-69=", $true, $false, $false, $false, $false, $true, 1, $false, "16+15=", 2) | Out-Null
$d.Content.Find.Execute("14+32=", $true, $false, $false, $false, $false, $true, 1, $false, "42-36=", 2) | Out-Null
$d.Content.Find.Execute("58-26=", $true, $false, $false, $false, $false, $true, 1, $false, "61-32=", 2) | Out-Null
$d.Content.Find.Execute("85-8=", $true, $false, $false, $false, $false, $true, 1, $false, "4+48=", 2) | Out-Null
$d.Content.Find.Execute("45-10=", $true, $false, $false, $false, $false, $true, 1, $false, "79-5=", 2) | Out-Null
$d.Content.Find.Execute("16+11=", $true, $false, $false, $false, $false, $true, 1, $false, "9-1=", 2) | Out-Null
$d.Content.Find.Execute("33+2=", $true, $false, $false, $false, $false, $true, 1, $false, "35+31=", 2) | Out-Null
$d.Content.Find.Execute("85+8=", $true, $false, $false, $false, $false, $true, 1, $false, "37+38=", 2) | Out-Null
$d.Content.Find.Execute("91-11=", $true, $false, $false, $false, $false, $true, 1, $false, "96-2=", 2) | Out-Null
$d.Content.Find.Execute("83-43=", $true, $false, $false, $false, $false, $true, 1, $false, "97-42=", 2) | Out-Null
$d.Content.Find.Execute("97-64=", $true, $false, $false, $false, $false, $true, 1, $false, "77-45=", 2) | Out-Null
